$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 1783
$wsExpo.Range("F6").Value = 149
$wsExpo.Range("F14").Value = 154
$wsExpo.Range("F15").Value = 24
$wsExpo.Range("F19").Value = 4938
$wsExpo.Range("G19").Value = 65
$wsExpo.Range("F20").Value = 49
$wsExpo.Range("F23").Value = 2237
$wsExpo.Range("F25").Value = 22
$wsExpo.Range("F26").Value = 2085

# --- Sheet "演出" ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 80

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1783
$wsAll.Range("F6").Value = 149
$wsAll.Range("F14").Value = 154
$wsAll.Range("F15").Value = 24
$wsAll.Range("F19").Value = 4939
$wsAll.Range("G19").Value = 65
$wsAll.Range("F20").Value = 80
$wsAll.Range("F21").Value = 49
$wsAll.Range("F25").Value = 2237
$wsAll.Range("F27").Value = 22
$wsAll.Range("F28").Value = 2085
